# Auto-generated edit script: updates cached price/profit values in the
# "Hades_Profits" leve-crafting sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect a refreshed Universalis market-data pull.

$wb = $excel.ActiveWorkbook


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3250.8713
$ws.Range("I15").Value = 3250.8713
$ws.Range("K15").Value = 9752.6139
$ws.Range("M15").Value = -9583.6139
$ws.Range("H80").Value = 623.73334
$ws.Range("I80").Value = 511.85715
$ws.Range("J80").Value = 721.625
$ws.Range("K80").Value = 1535.57145
$ws.Range("L80").Value = 2164.875
$ws.Range("M80").Value = -537.5714499999999
$ws.Range("N80").Value = -4160.875
$ws.Range("H83").Value = 623.73334
$ws.Range("I83").Value = 511.85715
$ws.Range("J83").Value = 721.625
$ws.Range("K83").Value = 4606.71435
$ws.Range("L83").Value = 6494.625
$ws.Range("M83").Value = 385.2856499999998
$ws.Range("N83").Value = -16478.625
$ws.Range("H107").Value = 500
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = ""
$ws.Range("N107").Value = -4340
$ws.Range("H132").Value = 733151.25
$ws.Range("I132").Value = 1835.0408
$ws.Range("J132").Value = 2723956.5
$ws.Range("K132").Value = 5505.1224
$ws.Range("L132").Value = 8171869.5
$ws.Range("M132").Value = -2975.1224
$ws.Range("N132").Value = -8176929.5
$ws.Range("H137").Value = 1787043.8
$ws.Range("I137").Value = 2381966.8
$ws.Range("K137").Value = 7145900.399999999
$ws.Range("M137").Value = -7143350.399999999
$ws.Range("H138").Value = 1962890.8
$ws.Range("I138").Value = 1426.0555
$ws.Range("J138").Value = 3403966.8
$ws.Range("K138").Value = 4278.166499999999
$ws.Range("L138").Value = 10211900.4
$ws.Range("M138").Value = 861.8335000000006
$ws.Range("N138").Value = -10222180.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 356.67
$ws.Range("I32").Value = 362.03226
$ws.Range("K32").Value = 362.03226
$ws.Range("M32").Value = -75.03226000000001
$ws.Range("H110").Value = 1699
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 1699
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 1699
$ws.Range("M110").Value = ""
$ws.Range("N110").Value = -5789
$ws.Range("H122").Value = 4275601
$ws.Range("I122").Value = 2225.8
$ws.Range("K122").Value = 6677.400000000001
$ws.Range("M122").Value = -4227.400000000001
$ws.Range("H132").Value = 34454.402
$ws.Range("I132").Value = 20726.92
$ws.Range("J132").Value = 91652.25
$ws.Range("K132").Value = 62180.75999999999
$ws.Range("L132").Value = 274956.75
$ws.Range("M132").Value = -59650.75999999999
$ws.Range("N132").Value = -280016.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 23731.111
$ws.Range("J138").Value = 23731.111
$ws.Range("L138").Value = 23731.111
$ws.Range("N138").Value = -34011.111

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 2600
$ws.Range("I32").Value = 2466.6667
$ws.Range("K32").Value = 2466.6667
$ws.Range("M32").Value = -2150.6667
$ws.Range("H58").Value = 20409652
$ws.Range("I58").Value = 24391658
$ws.Range("J58").Value = 1878.5
$ws.Range("K58").Value = 24391658
$ws.Range("L58").Value = 1878.5
$ws.Range("M58").Value = -24391455
$ws.Range("N58").Value = -2284.5
$ws.Range("H122").Value = 1782.1666
$ws.Range("I122").Value = 1670.0588
$ws.Range("J122").Value = 1928.7693
$ws.Range("K122").Value = 5010.1764
$ws.Range("L122").Value = 5786.3079
$ws.Range("M122").Value = -2560.1764
$ws.Range("N122").Value = -10686.3079
$ws.Range("H133").Value = 48975.066
$ws.Range("J133").Value = 49759
$ws.Range("L133").Value = 49759
$ws.Range("N133").Value = -54819
$ws.Range("H136").Value = 20409652
$ws.Range("I136").Value = 24391658
$ws.Range("J136").Value = 1878.5
$ws.Range("K136").Value = 73174974
$ws.Range("L136").Value = 5635.5
$ws.Range("M136").Value = -73172424
$ws.Range("N136").Value = -10735.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3395.1428
$ws.Range("I63").Value = 2141.6
$ws.Range("J63").Value = 4091.5557
$ws.Range("K63").Value = 6424.799999999999
$ws.Range("L63").Value = 12274.6671
$ws.Range("M63").Value = -5675.799999999999
$ws.Range("N63").Value = -13772.6671
$ws.Range("H64").Value = 3275.6206
$ws.Range("I64").Value = 1494.5
$ws.Range("J64").Value = 3407.5557
$ws.Range("K64").Value = 4483.5
$ws.Range("L64").Value = 10222.6671
$ws.Range("M64").Value = -4213.5
$ws.Range("N64").Value = -10762.6671
$ws.Range("H66").Value = 3395.1428
$ws.Range("I66").Value = 2141.6
$ws.Range("J66").Value = 4091.5557
$ws.Range("K66").Value = 19274.4
$ws.Range("L66").Value = 36824.0013
$ws.Range("M66").Value = -15530.4
$ws.Range("N66").Value = -44312.0013
$ws.Range("H67").Value = 3275.6206
$ws.Range("I67").Value = 1494.5
$ws.Range("J67").Value = 3407.5557
$ws.Range("K67").Value = 4483.5
$ws.Range("L67").Value = 10222.6671
$ws.Range("M67").Value = -3547.5
$ws.Range("N67").Value = -12094.6671
$ws.Range("H98").Value = 385.2
$ws.Range("I98").Value = 99.75
$ws.Range("K98").Value = 299.25
$ws.Range("M98").Value = 1198.75
$ws.Range("H107").Value = 761.91174
$ws.Range("I107").Value = 695.45
$ws.Range("J107").Value = 856.8570999999999
$ws.Range("K107").Value = 2086.35
$ws.Range("L107").Value = 2570.5713
$ws.Range("M107").Value = -166.3500000000004
$ws.Range("N107").Value = -6410.5713
$ws.Range("H121").Value = 40980830
$ws.Range("J121").Value = 46103250
$ws.Range("L121").Value = 138309750
$ws.Range("N121").Value = -138312370
$ws.Range("H129").Value = 3088310.8
$ws.Range("I129").Value = 1402
$ws.Range("J129").Value = 4904139.5
$ws.Range("K129").Value = 4206
$ws.Range("L129").Value = 14712418.5
$ws.Range("M129").Value = 794
$ws.Range("N129").Value = -14722418.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 657.6923
$ws.Range("I102").Value = 645.8333
$ws.Range("J102").Value = 800
$ws.Range("K102").Value = 645.8333
$ws.Range("L102").Value = 800
$ws.Range("M102").Value = 976.1667
$ws.Range("N102").Value = -4044

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 742.17145
$ws.Range("I22").Value = 491.53845
$ws.Range("J22").Value = 890.2727
$ws.Range("K22").Value = 491.53845
$ws.Range("L22").Value = 890.2727
$ws.Range("M22").Value = -196.53845
$ws.Range("N22").Value = -1480.2727
$ws.Range("H27").Value = 742.17145
$ws.Range("I27").Value = 491.53845
$ws.Range("J27").Value = 890.2727
$ws.Range("K27").Value = 491.53845
$ws.Range("L27").Value = 890.2727
$ws.Range("M27").Value = -384.53845
$ws.Range("N27").Value = -1104.2727
$ws.Range("H40").Value = 2290.4856
$ws.Range("I40").Value = 2202.0303
$ws.Range("J40").Value = 3750
$ws.Range("K40").Value = 2202.0303
$ws.Range("L40").Value = 3750
$ws.Range("M40").Value = -2066.0303
$ws.Range("N40").Value = -4022
$ws.Range("H132").Value = 43252
$ws.Range("I132").Value = 1436.1111
$ws.Range("J132").Value = 168699.67
$ws.Range("K132").Value = 4308.3333
$ws.Range("L132").Value = 506099.01
$ws.Range("M132").Value = -1778.3333
$ws.Range("N132").Value = -511159.01

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 35281.832
$ws.Range("I81").Value = 1546.1538
$ws.Range("J81").Value = 61079.707
$ws.Range("K81").Value = 3092.3076
$ws.Range("L81").Value = 122159.414
$ws.Range("M81").Value = -2031.3076
$ws.Range("N81").Value = -124281.414
$ws.Range("H84").Value = 35281.832
$ws.Range("I84").Value = 1546.1538
$ws.Range("J84").Value = 61079.707
$ws.Range("K84").Value = 15461.538
$ws.Range("L84").Value = 610797.0700000001
$ws.Range("M84").Value = -10157.538
$ws.Range("N84").Value = -621405.0700000001
